# "found another missing data item"
# Insert a new "Minimum Vapor Pressure" / "day-vap-pres-min" data item row
# into the Daily Vapor Pressure section (alphabetically between
# "Average Vapor Pressure" and the existing "Maximum Vapor Pressure" row),
# shifting all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Before: row18=Maximum Vapor Pressure, row19=Average Vapor Pressure.
# Insert a fresh row at 20 (shifts old rows 20.. down to 21..) and then
# re-write rows 18-20 in the new alphabetical order: Average, Maximum,
# Minimum.
$ws.Rows(20).Insert()

$ws.Range("A18").Value = "Average Vapor Pressure"
$ws.Range("B18").Value = "day-vap-pres-avg"
$ws.Range("C18").Value = "Daily"
$ws.Range("D18").Value = "WSN"

$ws.Range("A19").Value = "Maximum Vapor Pressure"
$ws.Range("B19").Value = "day-vap-pres-max"
$ws.Range("C19").Value = "Daily"
$ws.Range("D19").Value = "WSN"

$ws.Range("A20").Value = "Minimum Vapor Pressure"
$ws.Range("B20").Value = "day-vap-pres-min"
$ws.Range("C20").Value = "Daily"
$ws.Range("D20").Value = "WSN"

# Restore the view state recorded in the saved workbook.
[void]$ws.Range("B14").Select()
$excel.ActiveWindow.ScrollRow = 19
